# TC04_Canine_Filter_Diagnosis-Melanoma.xlsx
# "10 icdc scripts for jenkins" - update the FilesTab Cypher query (B4 on the
# "startup" sheet) to drop the `File Type` and `Breed` columns from the
# RETURN clause, and move the selection/scroll position down to row 4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$filesQuery = @'

MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
 MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
WHERE diag.disease_term IN ['Melanoma']
WITH DISTINCT f, parent, c, demo, diag, s
RETURN coalesce(f.file_name, '') AS `File Name`, 
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

# Defensive: make sure the text ends exactly on the `Study Code` backtick
# (no stray trailing newline), matching the source string.
$filesQuery = $filesQuery.TrimEnd("`r", "`n")

$ws.Range("B4").Value = $filesQuery

# Reflect the author's updated cursor/scroll position (row 4, column B).
$ws.Range("B4").Select()
